# Auto-generated Excel COM-interop script to apply the committed diff
# to Sheets/Asura_Profits.xlsx (the 8-sheet Asura Profits workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 1783.2858
$ws.Range("I5").Value = 1783.2858
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1783.2858
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1668.2858
$ws.Range("N5").Value = ""

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 771.5
$ws.Range("I18").Value = 771.5
$ws.Range("K18").Value = 771.5
$ws.Range("M18").Value = -487.5

# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 1798.1428
$ws.Range("I28").Value = 1931.1666
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 1931.1666
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -1446.1666
$ws.Range("N28").Value = -1970

# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 3415.818
$ws.Range("I32").Value = 5398.75
$ws.Range("J32").Value = 2282.7144
$ws.Range("K32").Value = 5398.75
$ws.Range("L32").Value = 2282.7144
$ws.Range("M32").Value = -5072.75
$ws.Range("N32").Value = -2934.7144

# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 18630896
$ws.Range("I70").Value = 41917492
$ws.Range("J70").Value = 1620
$ws.Range("K70").Value = 125752476
$ws.Range("L70").Value = 4860
$ws.Range("M70").Value = -125752206
$ws.Range("N70").Value = -5400

# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 18630896
$ws.Range("I73").Value = 41917492
$ws.Range("J73").Value = 1620
$ws.Range("K73").Value = 125752476
$ws.Range("L73").Value = 4860
$ws.Range("M73").Value = -125751540
$ws.Range("N73").Value = -6732

# Row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 1035.2084
$ws.Range("I127").Value = 569
$ws.Range("J127").Value = 1137.9323
$ws.Range("K127").Value = 1707
$ws.Range("L127").Value = 3413.7969
$ws.Range("M127").Value = 3253
$ws.Range("N127").Value = -13333.7969

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1380.2188
$ws.Range("I137").Value = 1254.56
$ws.Range("J137").Value = 1829
$ws.Range("K137").Value = 3763.68
$ws.Range("L137").Value = 5487
$ws.Range("M137").Value = -1213.68
$ws.Range("N137").Value = -10587

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2275.5596
$ws.Range("I138").Value = 1179.0454
$ws.Range("J138").Value = 3481.725
$ws.Range("K138").Value = 3537.1362
$ws.Range("L138").Value = 10445.175
$ws.Range("M138").Value = 1602.8638
$ws.Range("N138").Value = -20725.175

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 3157.4666
$ws.Range("I61").Value = 4852.4
$ws.Range("J61").Value = 2310
$ws.Range("K61").Value = 4852.4
$ws.Range("L61").Value = 2310
$ws.Range("M61").Value = -4640.4
$ws.Range("N61").Value = -2734

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1007.5
$ws.Range("I74").Value = 1033.6428
$ws.Range("J74").Value = 977
$ws.Range("K74").Value = 1033.6428
$ws.Range("L74").Value = 977
$ws.Range("M74").Value = -159.6428000000001
$ws.Range("N74").Value = -2725

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1007.5
$ws.Range("I77").Value = 1033.6428
$ws.Range("J77").Value = 977
$ws.Range("K77").Value = 5168.214
$ws.Range("L77").Value = 4885
$ws.Range("M77").Value = -800.2139999999999
$ws.Range("N77").Value = -13621

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 735.0606
$ws.Range("I97").Value = 580.9231
$ws.Range("J97").Value = 1307.5714
$ws.Range("K97").Value = 580.9231
$ws.Range("L97").Value = 1307.5714
$ws.Range("M97").Value = -84.92309999999998
$ws.Range("N97").Value = -2299.5714

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 7987.5454
$ws.Range("I132").Value = 11501.167
$ws.Range("J132").Value = 3771.2
$ws.Range("K132").Value = 34503.501
$ws.Range("L132").Value = 11313.6
$ws.Range("M132").Value = -31973.501
$ws.Range("N132").Value = -16373.6

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 3157.4666
$ws.Range("I136").Value = 4852.4
$ws.Range("J136").Value = 2310
$ws.Range("K136").Value = 14557.2
$ws.Range("L136").Value = 6930
$ws.Range("M136").Value = -12007.2
$ws.Range("N136").Value = -12030

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 127384.875
$ws.Range("I94").Value = 1050
$ws.Range("J94").Value = 169496.5
$ws.Range("K94").Value = 1050
$ws.Range("L94").Value = 169496.5
$ws.Range("M94").Value = -599
$ws.Range("N94").Value = -170398.5

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2248.775
$ws.Range("I134").Value = 1844.6
$ws.Range("J134").Value = 3461.3
$ws.Range("K134").Value = 5533.799999999999
$ws.Range("L134").Value = 10383.9
$ws.Range("M134").Value = -2998.799999999999
$ws.Range("N134").Value = -15453.9

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 4118526.5
$ws.Range("I58").Value = 5294248.5
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 5294248.5
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -5294045.5
$ws.Range("N58").Value = -3906

# Row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 32120
$ws.Range("J74").Value = 32120
$ws.Range("L74").Value = 32120
$ws.Range("N74").Value = -33868

# Row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 32120
$ws.Range("J77").Value = 32120
$ws.Range("L77").Value = 96360
$ws.Range("N77").Value = -105096

# Row 127 (Leve Item ID 35351)
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 4118526.5
$ws.Range("I136").Value = 5294248.5
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 15882745.5
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -15880195.5
$ws.Range("N136").Value = -15600

# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 40430.332
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 41734.125
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 41734.125
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -52094.125

$ws = $wb.Worksheets.Item("CUL")
# Row 87 (Leve Item ID 12864)
$ws.Range("H87").Value = 9500
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 9500
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 28500
$ws.Range("M87").Value = ""
$ws.Range("N87").Value = -30996

# Row 90 (Leve Item ID 12864)
$ws.Range("H90").Value = 9500
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 9500
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 85500
$ws.Range("M90").Value = ""
$ws.Range("N90").Value = -97980

# Row 104 (Leve Item ID 19807)
$ws.Range("H104").Value = 5929.091
$ws.Range("J104").Value = 5929.091
$ws.Range("L104").Value = 17787.273
$ws.Range("N104").Value = -23029.273

# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 453.85
$ws.Range("I113").Value = 417.9375
$ws.Range("J113").Value = 597.5
$ws.Range("K113").Value = 1253.8125
$ws.Range("L113").Value = 1792.5
$ws.Range("M113").Value = 916.1875
$ws.Range("N113").Value = -6132.5

# Row 120 (Leve Item ID 27877)
$ws.Range("H120").Value = 14691.333
$ws.Range("I120").Value = 10500
$ws.Range("J120").Value = 15529.6
$ws.Range("K120").Value = 31500
$ws.Range("L120").Value = 46588.8
$ws.Range("M120").Value = -26662
$ws.Range("N120").Value = -56264.8

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3061.5386
$ws.Range("I80").Value = 2987.5
$ws.Range("J80").Value = 3180
$ws.Range("K80").Value = 2987.5
$ws.Range("L80").Value = 3180
$ws.Range("M80").Value = -1989.5
$ws.Range("N80").Value = -5176

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3061.5386
$ws.Range("I83").Value = 2987.5
$ws.Range("J83").Value = 3180
$ws.Range("K83").Value = 14937.5
$ws.Range("L83").Value = 15900
$ws.Range("M83").Value = -9945.5
$ws.Range("N83").Value = -25884

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 5558
$ws.Range("J7").Value = 3990
$ws.Range("L7").Value = 3990
$ws.Range("N7").Value = -4214

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1716.6666
$ws.Range("I46").Value = 1650
$ws.Range("J46").Value = 1750
$ws.Range("K46").Value = 1650
$ws.Range("L46").Value = 1750
$ws.Range("M46").Value = -1462
$ws.Range("N46").Value = -2126

# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3200
$ws.Range("I68").Value = 3500
$ws.Range("J68").Value = 2900
$ws.Range("K68").Value = 3500
$ws.Range("L68").Value = 2900
$ws.Range("M68").Value = -2751
$ws.Range("N68").Value = -4398

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3200
$ws.Range("I71").Value = 3500
$ws.Range("J71").Value = 2900
$ws.Range("K71").Value = 17500
$ws.Range("L71").Value = 14500
$ws.Range("M71").Value = -13756
$ws.Range("N71").Value = -21988

# Row 110 (Leve Item ID 25809)
$ws.Range("H110").Value = 29362.666
$ws.Range("J110").Value = 29362.666
$ws.Range("L110").Value = 29362.666
$ws.Range("N110").Value = -37542.666

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 5558
$ws.Range("J126").Value = 3990
$ws.Range("L126").Value = 11970
$ws.Range("N126").Value = -16910

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 4333.3335
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376

# Row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 4333.3335
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 442.14285
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 465.55554
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 1396.66662
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -5236.66662

# Row 116 (Leve Item ID 26145)
$ws.Range("H116").Value = 35019.668
$ws.Range("J116").Value = 35019.668
$ws.Range("L116").Value = 35019.668
$ws.Range("N116").Value = -44197.668

# Row 117 (Leve Item ID 26162)
$ws.Range("H117").Value = 40204.5
$ws.Range("J117").Value = 40204.5
$ws.Range("L117").Value = 40204.5
$ws.Range("N117").Value = -49382.5

# Row 123 (Leve Item ID 34127)
$ws.Range("H123").Value = 24095.1
$ws.Range("J123").Value = 24095.1
$ws.Range("L123").Value = 24095.1
$ws.Range("N123").Value = -33895.1
